# Apply "changes to the data sheet to add new data" to the DataSheet.
# Rows 19-48 (A19:A48) of the DataSheet get replaced with 30 brand-new
# customer id strings. Rows 49-68 keep their existing values untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataSheet")

$newIds = @(
    "cus_Kv3XOsEZTFmjDW",
    "cus_Kv3XXdB1vGAo38",
    "cus_Kv3X020eO8lguC",
    "cus_Kv3XgOSFMBkpRh",
    "cus_Kv3XD7mvwnktxU",
    "cus_Kv3XRnnaf0jiP9",
    "cus_Kv3Xtns0fScCt4",
    "cus_Kv3XzDlmUMOz2P",
    "cus_Kv3XwRKGaq5AZG",
    "cus_Kv3Xg0y469rtXP",
    "cus_Kv3X8lCj4ZlsJJ",
    "cus_Kv3XWTseP4hkb6",
    "cus_Kv3XaKdkIuTmAu",
    "cus_Kv3XR8pASuv203",
    "cus_Kv3PU7F3UvzCLW",
    "cus_Kv3PQ4xeBzlevS",
    "cus_Kv3PHj455pR49r",
    "cus_Kv3PhbqE2WD89G",
    "cus_Kv3Pq0lSMhDpBs",
    "cus_Kv3PqN6r88hFcv",
    "cus_Kv3PKOCL9q68rf",
    "cus_Kv3PdNQrcKIaq9",
    "cus_Kv3P01RdpJZtFe",
    "cus_Kv3PS4M0ZmMXuZ",
    "cus_Kv3Pg8Ughzltui",
    "cus_Kv3PuHrHUDFOsK",
    "cus_Kv3PrtRFBZ4j3K",
    "cus_Kv3PjEKMVjD4eX",
    "cus_Kv3PDAtiKEbjDq",
    "cus_Kv3PYUrntTEKSr"
)

$startRow = 19
for ($i = 0; $i -lt $newIds.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newIds[$i]
}

$ws.Activate()
$ws.Range("A19:A48").Select()
